$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A54").Value = "Upload PDF file to the system as attach"
$ws.Range("B54").Value = "FAILED"
$ws.Range("C54").Value = "chrome"

$ws.Range("A55").Value = "Student should see PDF file to the system as attach"
$ws.Range("B55").Value = "FAILED"
$ws.Range("C55").Value = "chrome"

$ws.Range("A56").Value = "Upload PDF file to the system as attach"
$ws.Range("B56").Value = "PASSED"
$ws.Range("C56").Value = "chrome"

$ws.Range("A57").Value = "Student should see PDF file to the system as attach"
$ws.Range("B57").Value = "PASSED"
$ws.Range("C57").Value = "chrome"
